# F2_relay_banshee.xlsx — "Fixed issues with 81RF protective element"
# Changed default xls parameters to disable 81x protections.
#
# For relay setting rows 2-20 (sheet "relays"), the 81x (under/over
# frequency) protective-element columns are updated so the element is
# effectively disabled by default:
#   T / AH  (81x time dial / pickup #1) : 2   -> 100
#   U / AI  (81x pickup #2)             : 57  -> 10
#   V / AJ  (81x pickup #3)             : 0.01 -> 0.1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("relays")

for ($r = 2; $r -le 20; $r++) {
    $ws.Range("T$r").Value  = 100
    $ws.Range("U$r").Value  = 10
    $ws.Range("V$r").Value  = 0.1

    $ws.Range("AH$r").Value = 100
    $ws.Range("AI$r").Value = 10
    $ws.Range("AJ$r").Value = 0.1
}

# Minor row-height touch-up that accompanied the resave (rows slightly
# taller to match the updated default grid metrics).
$ws.Rows.Item(1).RowHeight = 43.15
$ws.Range("2:21").RowHeight = 14.45

# Leave the view scrolled/selected over the newly-edited 81x block, as in
# the source workbook after the edit.
$ws.Range("AH2:AJ20").Select()
